$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 12.85
$ws.Range("E7").Value = 13.008
$ws.Range("D8").Value = -7.887
$ws.Range("B12").Value = 5.568000000000001
$ws.Range("D12").Value = -8.261000000000001
$ws.Range("D14").Value = -8.263
$ws.Range("E19").Value = 12.452
$ws.Range("E21").Value = 12.923
$ws.Range("D22").Value = -8.106000000000002
$ws.Range("E24").Value = 12.848
